$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '67.465.16'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.626.87'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '594.15'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '168.28'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +1.03%  '
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -2.28%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.627.57'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -1.79%  '
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -1.86%  '
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.19%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.364'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +1.57%  '
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.22'
$c.ClearFormats()
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.ClearFormats()
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.63'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -0.71%  '
$c.ClearFormats()
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.107.19'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -1.79%  '
$c.ClearFormats()
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.0000182'
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -1.31%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '67.371.57'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.627.67'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -1.66%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '12.02'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +2.38%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '8.04'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +4.45%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '357.20'
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.32'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -1.31%  '
$c.ClearFormats()
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '4.67'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -2.95%  '
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -4.06%  '
$c.ClearFormats()
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '10.28'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +2.94%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '69.64'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -1.67%  '
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -1.60%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '546.18'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -2.12%  '
$c.ClearFormats()
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.92'
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -1.27%  '
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -2.88%  '
$c.ClearFormats()
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.89'
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -1.93%  '
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +4.48%  '
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.10%  '
$c.ClearFormats()
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.50'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -3.19%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '156.66'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +1.22%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '19.01'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -2.71%  '
$c.ClearFormats()
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.366'
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -1.84%  '
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.77%  '
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +1.89%  '
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -1.33%  '
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -3.62%  '
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.ClearFormats()
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '152.67'
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.579'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -1.98%  '
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '3.79'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -1.28%  '
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -1.33%  '
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.20%  '
$c.ClearFormats()
